$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3600212038047717
$ws.Range("C2").Value = 0.4297948612912168
$ws.Range("D2").Value = 0.3689482115053426
$ws.Range("E2").Value = 0.6074110729196025
$ws.Range("F2").Value = 0.507684125024779

$ws.Range("B3").Value = 0.1144269351476378
$ws.Range("C3").Value = 0.2119524186144291
$ws.Range("D3").Value = 0.09979453641461808
$ws.Range("E3").Value = 0.3159027325216072
$ws.Range("F3").Value = 0.3103779289964921

$ws.Range("B4").Value = 0.04590311845355863
$ws.Range("C4").Value = 0.2393856741142939
$ws.Range("D4").Value = 0.09493667408959593
$ws.Range("E4").Value = 0.3081179548315806
$ws.Range("F4").Value = 0.3337596341186294

$ws.Range("B5").Value = 0.01339574438561458
$ws.Range("C5").Value = 0.07924522757218994
$ws.Range("D5").Value = 0.006459252060612895
$ws.Range("E5").Value = 0.08036947219319594
$ws.Range("F5").Value = 0.1120696755859333
